$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (LC0049 - GroupAnagrams notebook)
$ws.Range("A19").Value = "x"
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = "Array"
$ws.Range("H19").Value = "Hash Table"
$ws.Range("I19").Value = "String"

# Row 20 (LC0049 - group-anagrams.py)
$ws.Range("A20").Value = "x"
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = "Array"
$ws.Range("H20").Value = "Hash Table"
$ws.Range("I20").Value = "String"

# Row 29 (LC0073 - SetMatrixZeroes notebook)
$ws.Range("A29").Value = "x"
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = "Hash Table"
$ws.Range("I29").Value = "Matrix"

# Writing cell values triggers an auto-fit recalculation of row height in
# this runtime; restore the original (unchanged-per-diff) row heights.
$ws.Rows.Item(19).RowHeight = 15.95
$ws.Rows.Item(20).RowHeight = 15.95
$ws.Rows.Item(29).RowHeight = 15.95

Write-Output "done"
